# New Submission Synced: 2026-02-09 19:24:17
# A new Google-Forms-style response row was added to the "JSS 3D" results
# sheet: Timestamp, Full Name, Admission No (kept as text, e.g. "1"), AI Score.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JSS 3D")

# Force column C's new cell to stay text (the Admission No for this
# submission is the numeric-looking string "1", same as row 3/4/5 which
# already store their admission numbers as text) without leaving a
# lingering custom number format on the cell.
$ws.Range("C7").NumberFormat = "@"

$ws.Cells.Item(7, 1).Value = "2026-02-09 19:24:17"
$ws.Cells.Item(7, 2).Value = "Mustapha Ali abbatar "
$ws.Cells.Item(7, 3).Value = "1"
$ws.Cells.Item(7, 4).Value = 10

# Reset the style reference back to the default "Normal" style now that the
# value has been committed as text, so no stray style id is left on C7.
$ws.Range("C7").Style = "Normal"
